# "Move or Copy" -> "Create a copy" of the Slovakia sheet, dropped at the
# end of the tab strip, then renamed to "Italy" and its two market-specific
# cells (B2 = market name, B4 = Jira/user-story ref) updated for Italy.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Slovakia")
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$italy = $wb.ActiveSheet
$italy.Name = "Italy"

$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2237"

# Mirror the selection state Excel leaves behind after a sheet copy: the
# source tab ends up with its whole grid selected, while the freshly
# created (and now active) copy keeps the focus on B4.
$template.Activate()
$template.Cells.Select()

$italy.Activate()
$italy.Range("B4").Select()
